# Apply color-wise numeric "Values" to each sheet's column B (rows 2-15),
# replacing the placeholder color-name text that was there before.
# The "Blue" sheet instead has its B2:B15 cells cleared entirely (no data).

$wb = $excel.ActiveWorkbook

# Sheet "D Green"
$ws = $wb.Worksheets.Item("D Green")
$values = @(0, 0, 0, 21, 0, 10, 0, 17, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $values[$i]
}

# Sheet "Green"
$ws = $wb.Worksheets.Item("Green")
$values = @(9, 29, 97, 64, 55, 81, 147, 85, 50, 79, 95, 51, 72, 92)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $values[$i]
}

# Sheet "Yellow"
$ws = $wb.Worksheets.Item("Yellow")
$values = @(77.25, 55.5, 36.75, 59.25, 62.25, 77.25, 25.5, 41.25, 60, 60.75, 56.25, 76.5, 48.75, 44.25)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $values[$i]
}

# Sheet "Orange"
$ws = $wb.Worksheets.Item("Orange")
$values = @(32.5, 28, 11, 9.5, 12.5, 4, 7, 13.5, 14.5, 13.5, 13, 18.5, 23, 17)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $values[$i]
}

# Sheet "Brown"
$ws = $wb.Worksheets.Item("Brown")
$values = @(2.5, 6.5, 2.5, 4, 8.5, 1.5, 2.5, 3, 5.75, 2.25, 3, 4.25, 4.5, 5.5)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $values[$i]
}

# Sheet "Red" (all zeros)
$ws = $wb.Worksheets.Item("Red")
$ws.Range("B2:B15").Value = 0

# Sheet "Default Red" (all zeros)
$ws = $wb.Worksheets.Item("Default Red")
$ws.Range("B2:B15").Value = 0

# Sheet "Blue": remove the B2:B15 cell contents entirely
$ws = $wb.Worksheets.Item("Blue")
$ws.Range("B2:B15").ClearContents()
